# Updated Process to open the websites with the data files
#
# Replace the placeholder "website" rows with the real NHANES dataset
# pages (one row per NHANES data component), update the header/label
# text, and give the first data-page link the "pasted from the web"
# look (Open Sans 9pt / #212529) that Excel applies when a hyperlink
# is pasted in from a browser.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 : Demographic ---------------------------------------------
$ws.Range("A2").Value = "https://wwwn.cdc.gov/nchs/nhanes/search/datapage.aspx?Component=Demographics&CycleBeginYear=2017"
$ws.Range("B2").Value = "Demographic Data"
$ws.Range("B2").Font.Name = "Open Sans"
$ws.Range("B2").Font.Size = 9
$ws.Range("B2").Font.Color = 2696481

# --- Row 3 : Dietary ---------------------------------------------------
$ws.Range("A3").Value = "https://wwwn.cdc.gov/nchs/nhanes/search/datapage.aspx?Component=Dietary&CycleBeginYear=2017"
$ws.Range("B3").Value = "Dietary Data"

# --- Row 4 : Examination ------------------------------------------------
$ws.Range("A4").Value = "https://wwwn.cdc.gov/nchs/nhanes/search/datapage.aspx?Component=Examination&CycleBeginYear=2017"
$ws.Range("B4").Value = "Examination Data"

# --- Row 5 : Laboratory (new row) ---------------------------------------
$ws.Range("A5").Value = "https://wwwn.cdc.gov/nchs/nhanes/search/datapage.aspx?Component=Laboratory&CycleBeginYear=2017"
$ws.Range("B5").Value = "Laboratory Data"

# --- Row 6 : Questionnaire (new row) -------------------------------------
$ws.Range("A6").Value = "https://wwwn.cdc.gov/nchs/nhanes/search/datapage.aspx?Component=Questionnaire&CycleBeginYear=2017"
$ws.Range("B6").Value = "Questionnaire Data"

# --- Row 7 : Limited Access Data (new row) -------------------------------
$ws.Range("A7").Value = "https://wwwn.cdc.gov/nchs/nhanes/search/datapage.aspx?Component=LimitedAccess&CycleBeginYear=2017"
$ws.Range("B7").Value = "Limited Access Data"

# Final selection left on A3, matching the saved workbook state.
$null = $ws.Range("A3").Select()
